$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data grid (header + 2 rows), columns A-F:
#   A=url, B=topic_id, C=topic, D=question, E=level, F=platform
$data = @(
    @("url", "topic_id", "topic", "question", "level", "platform"),
    @("abcd", "fc1c3f36164311eea88ae3300d621ca4", "arrays", "asdf", "easy", "codechef"),
    @("adcde", "55c324b8164511eea88ae3300d621ca4", "twoPointers", "sadd", "easy", "codechef")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$wb.Save()
